$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.845.13'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '2.954.07'
$ws.Range('E3').Value = '  +2.32%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '555.05'
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('D6').Value = '133.73'
$ws.Range('E6').Value = '  +10.29%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +4.14%  '
$ws.Range('D9').Value = '2.945.25'
$ws.Range('E9').Value = '  +2.33%  '
$ws.Range('E10').Value = '  +4.04%  '
$ws.Range('E11').Value = '  +2.25%  '
$ws.Range('E12').Value = '  +3.91%  '
$ws.Range('E13').Value = '  +5.12%  '
$ws.Range('D14').Value = '32.83'
$ws.Range('E14').Value = '  +4.65%  '
$ws.Range('E15').Value = '  +3.02%  '
$ws.Range('D16').Value = '3.438.12'
$ws.Range('E16').Value = '  +2.21%  '
$ws.Range('D17').Value = '6.98'
$ws.Range('E17').Value = '  +7.49%  '
$ws.Range('D18').Value = '2.948.74'
$ws.Range('E18').Value = '  +2.36%  '
$ws.Range('D19').Value = '57.807.85'
$ws.Range('E19').Value = '  +0.82%  '
$ws.Range('D20').Value = '417.73'
$ws.Range('E20').Value = '  +2.32%  '
$ws.Range('D21').Value = '13.43'
$ws.Range('E21').Value = '  +5.30%  '
$ws.Range('E22').Value = '  +7.46%  '
$ws.Range('D23').Value = '13.43'
$ws.Range('E23').Value = '  +7.12%  '
$ws.Range('E24').Value = '  +4.80%  '
$ws.Range('D25').Value = '79.70'
$ws.Range('E25').Value = '  +4.03%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('E28').Value = '  +1.86%  '
$ws.Range('E29').Value = '  +7.17%  '
$ws.Range('E30').Value = '  +6.27%  '
$ws.Range('E31').Value = '  +4.10%  '
$ws.Range('E32').Value = '  -0.81%  '
$ws.Range('D33').Value = '0.0969'
$ws.Range('E33').Value = '  +3.33%  '
$ws.Range('E34').Value = '  +7.14%  '
$ws.Range('D35').Value = '0.953'
$ws.Range('E35').Value = '  +7.16%  '
$ws.Range('E36').Value = '  +2.47%  '
$ws.Range('D37').Value = '0.0₃0702'
$ws.Range('E37').Value = '  +14.33%  '
$ws.Range('D38').Value = '8.92'
$ws.Range('E38').Value = '  +7.50%  '
$ws.Range('D39').Value = '48.14'
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('D40').Value = '2.72'
$ws.Range('E40').Value = '  +16.83%  '
$ws.Range('D41').Value = '385.32'
$ws.Range('E41').Value = '  +6.56%  '
$ws.Range('E42').Value = '  +3.07%  '
$ws.Range('E43').Value = '  +1.51%  '
$ws.Range('D44').Value = '2.717.18'
$ws.Range('E44').Value = '  +4.62%  '
$ws.Range('D46').Value = '124.98'
$ws.Range('E46').Value = '  +5.71%  '
$ws.Range('D47').Value = '0.237'
$ws.Range('E47').Value = '  +4.38%  '
$ws.Range('D48').Value = '1.99'
$ws.Range('E48').Value = '  +4.08%  '
$ws.Range('E49').Value = '  +1.87%  '
$ws.Range('D50').Value = '22.89'
$ws.Range('E50').Value = '  +3.55%  '
$ws.Range('E51').Value = '  +3.75%  '
